# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) on each crafting-class sheet with freshly pulled Universalis prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: Distill, My Heart
$ws.Range("H9").Value = 73.5
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# Row 12: Don't Be So Tallow
$ws.Range("H12").Value = 199.75
$ws.Range("I12").Value = 199.66667
$ws.Range("K12").Value = 199.66667
$ws.Range("M12").Value = -29.66667000000001

# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 157.3
$ws.Range("I15").Value = 157.3
$ws.Range("K15").Value = 471.9
$ws.Range("M15").Value = -302.9

# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 3629.6
$ws.Range("I74").Value = 3287.25
$ws.Range("J74").Value = 4999
$ws.Range("K74").Value = 3287.25
$ws.Range("L74").Value = 4999
$ws.Range("M74").Value = -2351.25
$ws.Range("N74").Value = -6871

# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 3629.6
$ws.Range("I77").Value = 3287.25
$ws.Range("J77").Value = 4999
$ws.Range("K77").Value = 16436.25
$ws.Range("L77").Value = 24995
$ws.Range("M77").Value = -11756.25
$ws.Range("N77").Value = -34355

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 57687.332
$ws.Range("J87").Value = 57687.332
$ws.Range("L87").Value = 57687.332
$ws.Range("N87").Value = -60183.332

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 57687.332
$ws.Range("J90").Value = 57687.332
$ws.Range("L90").Value = 173061.996
$ws.Range("N90").Value = -185541.996

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 1210.44
$ws.Range("I92").Value = 1215.3
$ws.Range("J92").Value = 1191
$ws.Range("K92").Value = 1215.3
$ws.Range("L92").Value = 1191
$ws.Range("M92").Value = 32.70000000000005
$ws.Range("N92").Value = -3687

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 4904982.5
$ws.Range("I132").Value = 2804.6597
$ws.Range("J132").Value = 62505576
$ws.Range("K132").Value = 8413.9791
$ws.Range("L132").Value = 187516728
$ws.Range("M132").Value = -5883.9791
$ws.Range("N132").Value = -187521788

# Row 138: All-night Crafting
$ws.Range("H138").Value = 4633081.5
$ws.Range("I138").Value = 10103384
$ws.Range("J138").Value = 4364.4873
$ws.Range("K138").Value = 30310152
$ws.Range("L138").Value = 13093.4619
$ws.Range("M138").Value = -30305012
$ws.Range("N138").Value = -23373.4619

# Row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 44978.09
$ws.Range("J139").Value = 44978.09
$ws.Range("L139").Value = 44978.09
$ws.Range("N139").Value = -55258.09

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 19180.945
$ws.Range("I32").Value = 18749.076
$ws.Range("J32").Value = 26666.666
$ws.Range("K32").Value = 18749.076
$ws.Range("L32").Value = 26666.666
$ws.Range("M32").Value = -18462.076
$ws.Range("N32").Value = -27240.666

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 3001.6
$ws.Range("I132").Value = 2331.0588
$ws.Range("J132").Value = 4426.5
$ws.Range("K132").Value = 6993.176399999999
$ws.Range("L132").Value = 13279.5
$ws.Range("M132").Value = -4463.176399999999
$ws.Range("N132").Value = -18339.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 894.8
$ws.Range("I94").Value = 707.4286
$ws.Range("J94").Value = 1332
$ws.Range("K94").Value = 707.4286
$ws.Range("L94").Value = 1332
$ws.Range("M94").Value = -256.4286
$ws.Range("N94").Value = -2234

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 3205.2698
$ws.Range("I31").Value = 2028.7368
$ws.Range("J31").Value = 4993.6
$ws.Range("K31").Value = 2028.7368
$ws.Range("L31").Value = 4993.6
$ws.Range("M31").Value = -1733.7368
$ws.Range("N31").Value = -5583.6

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 3205.2698
$ws.Range("I34").Value = 2028.7368
$ws.Range("J34").Value = 4993.6
$ws.Range("K34").Value = 2028.7368
$ws.Range("L34").Value = 4993.6
$ws.Range("M34").Value = -1826.7368
$ws.Range("N34").Value = -5397.6

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3406.6086
$ws.Range("I58").Value = 979.3333
$ws.Range("J58").Value = 6054.5454
$ws.Range("K58").Value = 979.3333
$ws.Range("L58").Value = 6054.5454
$ws.Range("M58").Value = -776.3333
$ws.Range("N58").Value = -6460.5454

# Row 99: O Pine
$ws.Range("H99").Value = 2229.75
$ws.Range("I99").Value = 1569.4
$ws.Range("J99").Value = 2890.1
$ws.Range("K99").Value = 1569.4
$ws.Range("L99").Value = 2890.1
$ws.Range("M99").Value = -71.40000000000009
$ws.Range("N99").Value = -5886.1

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 2119.5
$ws.Range("I122").Value = 2093.4
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 6280.200000000001
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -3830.200000000001
$ws.Range("N122").Value = -11650

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2229.75
$ws.Range("I126").Value = 1569.4
$ws.Range("J126").Value = 2890.1
$ws.Range("K126").Value = 4708.200000000001
$ws.Range("L126").Value = 8670.299999999999
$ws.Range("M126").Value = -2238.200000000001
$ws.Range("N126").Value = -13610.3

# Row 136: Turali Quality
$ws.Range("H136").Value = 3406.6086
$ws.Range("I136").Value = 979.3333
$ws.Range("J136").Value = 6054.5454
$ws.Range("K136").Value = 2937.9999
$ws.Range("L136").Value = 18163.6362
$ws.Range("M136").Value = -387.9998999999998
$ws.Range("N136").Value = -23263.6362

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1655.8125
$ws.Range("I97").Value = 1126.25
$ws.Range("K97").Value = 1126.25
$ws.Range("M97").Value = -630.25

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 2294.7046
$ws.Range("I102").Value = 2358.4062
$ws.Range("K102").Value = 2358.4062
$ws.Range("M102").Value = -736.4061999999999

# Row 132: On Board for Lar
$ws.Range("H132").Value = 12233.538
$ws.Range("I132").Value = 19360.572
$ws.Range("J132").Value = 3918.6667
$ws.Range("K132").Value = 58081.716
$ws.Range("L132").Value = 11756.0001
$ws.Range("M132").Value = -55551.716
$ws.Range("N132").Value = -16816.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 2475.0557
$ws.Range("I16").Value = 2361.6155
$ws.Range("J16").Value = 2770
$ws.Range("K16").Value = 2361.6155
$ws.Range("L16").Value = 2770
$ws.Range("M16").Value = -2191.6155
$ws.Range("N16").Value = -3110

# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1551.2307
$ws.Range("I22").Value = 1100
$ws.Range("J22").Value = 1633.2727
$ws.Range("K22").Value = 1100
$ws.Range("L22").Value = 1633.2727
$ws.Range("M22").Value = -805
$ws.Range("N22").Value = -2223.2727

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1551.2307
$ws.Range("I27").Value = 1100
$ws.Range("J27").Value = 1633.2727
$ws.Range("K27").Value = 1100
$ws.Range("L27").Value = 1633.2727
$ws.Range("M27").Value = -993
$ws.Range("N27").Value = -1847.2727

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 2081.1
$ws.Range("I68").Value = 1876
$ws.Range("K68").Value = 1876
$ws.Range("M68").Value = -1127

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 2081.1
$ws.Range("I71").Value = 1876
$ws.Range("K71").Value = 9380
$ws.Range("M71").Value = -5636

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 1643.9166
$ws.Range("I82").Value = 1530.5834
$ws.Range("J82").Value = 1757.25
$ws.Range("K82").Value = 1530.5834
$ws.Range("L82").Value = 1757.25
$ws.Range("M82").Value = -1169.5834
$ws.Range("N82").Value = -2479.25

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 1643.9166
$ws.Range("I85").Value = 1530.5834
$ws.Range("J85").Value = 1757.25
$ws.Range("K85").Value = 1530.5834
$ws.Range("L85").Value = 1757.25
$ws.Range("M85").Value = -282.5834
$ws.Range("N85").Value = -4253.25

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 10006470
$ws.Range("I132").Value = 5285.517
$ws.Range("J132").Value = 23817628
$ws.Range("K132").Value = 15856.551
$ws.Range("L132").Value = 71452884
$ws.Range("M132").Value = -13326.551
$ws.Range("N132").Value = -71457944

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 1704.0222
$ws.Range("I136").Value = 929.6
$ws.Range("J136").Value = 4414.5
$ws.Range("K136").Value = 2788.8
$ws.Range("L136").Value = 13243.5
$ws.Range("M136").Value = -238.8000000000002
$ws.Range("N136").Value = -18343.5

Write-Host "Ultima Profits sheets updated."
